$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Aoking UNISEX - Rucksack - mint" item in row 29 with the
# new "Aoking Tagesrucksack - mint" item (new image + link, same price).
$ws.Range("A29").Value = "Aoking Tagesrucksack - mint"
$ws.Range("B29").Value = "https://img01.ztat.net/article/spp-media-p1/5b94485402654a06821e70f335da2a99/8c5feb95558746c2b5143245f8749a1d.jpg?imwidth=1800"
$ws.Range("C29").Value = "https://www.zalando.ch/aoking-tagesrucksack-mint-ao054h00o-m11.html"
$ws.Range("D29").Value = "75 CHF"

# Append two new wishlist items as new rows 40 and 41. Column A of the
# "category header" rows (e.g. A38) uses a bold + orange-fill style, so
# copy that formatting onto the new A40/A41 cells.
$ws.Range("A38").Copy()
$ws.Range("A40").PasteSpecial(-4122)
$ws.Range("A41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A40").Value = "CALVIN KLEIN ETERNITY FOR HER EAU DE PARFUM"
$ws.Range("B40").Value = "https://img01.ztat.net/article/spp-media-p1/82b5fb685c434f4cacede01af31d3d35/9f21bdd22d9f44938725464b8b8524c2.jpg?imwidth=1800&filter=packshot"
$ws.Range("C40").Value = "https://www.zalando.ch/calvin-klein-fragrances-calvin-klein-eternity-for-her-eau-de-parfum-eau-de-parfum-c4p31i005-s11.html"
$ws.Range("D40").Value = "42 CHF"

$ws.Range("A41").Value = "Givenchy Ange ou Démon (Etrange)"
$ws.Range("B41").Value = "https://static01.galaxus.com/productimages/4/9/4/3/1/8/8/1/8/0/4/5/6/0/8/6/2/3/2/01993825-56b5-778a-b12e-e3be33af0bf6_2880.avif"
$ws.Range("C41").Value = "https://www.galaxus.ch/en/s6/product/givenchy-ange-ou-dmon-etrange-eau-de-parfum-100-ml-fragrances-22870522"
$ws.Range("D41").Value = "65 CHF"

# Update the sheet selection / active cell to match the edited workbook.
$ws.Range("A30").Select()
